# Fixing test data: remove the "Control" row from the "On-Site" sheet
# and delete the entire "Replacement" sheet.

$wb = $excel.ActiveWorkbook

# --- On-Site sheet: delete the row containing "Control" ---
$onSite = $wb.Worksheets.Item("On-Site")
$onSite.Activate()

$controlRow = $null
for ($r = 1; $r -le $onSite.UsedRange.Rows.Count; $r++) {
    if ($onSite.Cells.Item($r, 1).Value2 -eq "Control") {
        $controlRow = $r
        break
    }
}

if ($controlRow -ne $null) {
    $onSite.Rows.Item($controlRow).Delete() | Out-Null
}

# Select the row that now occupies the old "Control" row position
$onSite.Range("A" + $controlRow + ":XFD" + $controlRow).Select() | Out-Null

# --- Delete the "Replacement" sheet entirely ---
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Replacement").Delete() | Out-Null
$excel.DisplayAlerts = $true
